$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp shown at the top of the sheet.
$ws.Range("A1").Value = "Datos actualizados a 10 de Octubre de 2020 a las 05:49"

# Belgica moves up in the ranking (new, larger totals) ahead of Rumania,
# Marruecos and Ecuador, which are each pushed one position down. This
# mirrors a re-sort of the country list by total cases, so rows 31-34
# need to be rewritten with the post-sort country name + stats.
$ws.Cells.Item(31, 1).Value = "Belgica"
$ws.Cells.Item(31, 2).Value = 148981
$ws.Cells.Item(31, 3).Value = 5385
$ws.Cells.Item(31, 4).Value = 20072
$ws.Cells.Item(31, 5).Value = 118758
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(31, 7).Value = 25
$ws.Cells.Item(31, 8).Value = 10151

$ws.Cells.Item(32, 1).Value = "Rumania"
$ws.Cells.Item(32, 2).Value = 148886
$ws.Cells.Item(32, 3).Value = 0
$ws.Cells.Item(32, 4).Value = 114792
$ws.Cells.Item(32, 5).Value = 28795
$ws.Cells.Item(32, 6).Value = 0
$ws.Cells.Item(32, 7).Value = 0
$ws.Cells.Item(32, 8).Value = 5299

$ws.Cells.Item(33, 1).Value = "Marruecos"
$ws.Cells.Item(33, 2).Value = 146398
$ws.Cells.Item(33, 3).Value = 0
$ws.Cells.Item(33, 4).Value = 123022
$ws.Cells.Item(33, 5).Value = 20846
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(33, 7).Value = 0
$ws.Cells.Item(33, 8).Value = 2530

$ws.Cells.Item(34, 1).Value = "Ecuador"
$ws.Cells.Item(34, 2).Value = 145848
$ws.Cells.Item(34, 3).Value = 0
$ws.Cells.Item(34, 4).Value = 120511
$ws.Cells.Item(34, 5).Value = 13162
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 12175

# Belice moves up ahead of Benin and Guinea-Bisau, which are each pushed
# one position down.
$ws.Cells.Item(153, 1).Value = "Belice"
$ws.Cells.Item(153, 2).Value = 2427
$ws.Cells.Item(153, 3).Value = 54
$ws.Cells.Item(153, 4).Value = 1487
$ws.Cells.Item(153, 5).Value = 905
$ws.Cells.Item(153, 6).Value = 0
$ws.Cells.Item(153, 7).Value = 1
$ws.Cells.Item(153, 8).Value = 35

$ws.Cells.Item(154, 1).Value = "Benin"
$ws.Cells.Item(154, 2).Value = 2411
$ws.Cells.Item(154, 3).Value = 0
$ws.Cells.Item(154, 4).Value = 1973
$ws.Cells.Item(154, 5).Value = 397
$ws.Cells.Item(154, 6).Value = 0
$ws.Cells.Item(154, 7).Value = 0
$ws.Cells.Item(154, 8).Value = 41

$ws.Cells.Item(155, 1).Value = "Guinea-Bisau"
$ws.Cells.Item(155, 2).Value = 2385
$ws.Cells.Item(155, 3).Value = 0
$ws.Cells.Item(155, 4).Value = 1728
$ws.Cells.Item(155, 5).Value = 617
$ws.Cells.Item(155, 6).Value = 0
$ws.Cells.Item(155, 7).Value = 0
$ws.Cells.Item(155, 8).Value = 40

# Updated stats for San Martin (Parte Holandesa); no reordering involved.
$ws.Cells.Item(172, 2).Value = 703
$ws.Cells.Item(172, 3).Value = 4
$ws.Cells.Item(172, 4).Value = 622
$ws.Cells.Item(172, 5).Value = 59

# Updated stats for Mongolia; no reordering involved.
$ws.Cells.Item(186, 4).Value = 310
$ws.Cells.Item(186, 5).Value = 5
